# Auto-generated edit script applying the cryptos.xlsx price/volume diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.335.79"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -3.22%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = "'1.772.80"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -1.99%  "
$ws.Range("E3").ClearFormats()

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E4").ClearFormats()

$ws.Range("E5").Value = "'  +0.12%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D6").Value = "'305.89"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -1.37%  "
$ws.Range("E6").ClearFormats()

$ws.Range("D7").Value = "'0.4226"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +0.50%  "
$ws.Range("E7").ClearFormats()

$ws.Range("D8").Value = "'0.3596"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +0.93%  "
$ws.Range("E8").ClearFormats()

$ws.Range("E9").Value = "'  -0.05%  "
$ws.Range("E9").ClearFormats()

$ws.Range("D10").Value = "'0.8353"
$ws.Range("D10").ClearFormats()

$ws.Range("E11").Value = "'  +0.82%  "
$ws.Range("E11").ClearFormats()

$ws.Range("D12").Value = "'1.784.08"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -1.79%  "
$ws.Range("E12").ClearFormats()

$ws.Range("D13").Value = "'6.440"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +1.10%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").Value = "'5.226"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -1.64%  "
$ws.Range("E14").ClearFormats()

$ws.Range("D15").Value = "'0.06860"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +0.12%  "
$ws.Range("E16").ClearFormats()

$ws.Range("D17").Value = "'78.72"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -2.84%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").Value = "'0.000008613"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -1.64%  "
$ws.Range("E18").ClearFormats()

$ws.Range("E19").Value = "'  +0.11%  "
$ws.Range("E19").ClearFormats()

$ws.Range("D20").Value = "'14.88"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -1.61%  "
$ws.Range("E20").ClearFormats()

$ws.Range("D21").Value = "'26.387.79"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -3.10%  "
$ws.Range("E21").ClearFormats()

$ws.Range("D22").Value = "'5.070"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -0.71%  "
$ws.Range("E22").ClearFormats()

$ws.Range("D23").Value = "'10.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +0.95%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").Value = "'2.016.87"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -0.73%  "
$ws.Range("E24").ClearFormats()

$ws.Range("D25").Value = "'152.30"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.01%  "
$ws.Range("E25").ClearFormats()

$ws.Range("D26").Value = "'1.795"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -8.84%  "
$ws.Range("E26").ClearFormats()

$ws.Range("D27").Value = "'17.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -1.01%  "
$ws.Range("E27").ClearFormats()

$ws.Range("D28").Value = "'5.054"
$ws.Range("D28").ClearFormats()

$ws.Range("D29").Value = "'114.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +0.79%  "
$ws.Range("E29").ClearFormats()

$ws.Range("D30").Value = "'1.822"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +8.45%  "
$ws.Range("E30").ClearFormats()

$ws.Range("D31").Value = "'0.08842"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -0.65%  "
$ws.Range("E31").ClearFormats()

$ws.Range("D32").Value = "'0.7244"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -1.63%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").Value = "'1.112"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +0.72%  "
$ws.Range("E33").ClearFormats()

$ws.Range("D34").Value = "'4.315"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -2.48%  "
$ws.Range("E34").ClearFormats()

$ws.Range("D35").Value = "'1.002"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +0.13%  "
$ws.Range("E35").ClearFormats()

$ws.Range("D36").Value = "'2.736"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -7.19%  "
$ws.Range("E36").ClearFormats()

$ws.Range("D37").Value = "'1.093"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +0.94%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").Value = "'0.05128"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -0.87%  "
$ws.Range("E38").ClearFormats()

$ws.Range("D39").Value = "'0.01879"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -1.33%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").Value = "'0.1607"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").Value = "'0.4895"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -1.54%  "
$ws.Range("E41").ClearFormats()

$ws.Range("D42").Value = "'2.618"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -3.51%  "
$ws.Range("E42").ClearFormats()

$ws.Range("D43").Value = "'6.375"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +1.85%  "
$ws.Range("E43").ClearFormats()

$ws.Range("D44").Value = "'7.951"
$ws.Range("D44").ClearFormats()

$ws.Range("D45").Value = "'104.55"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -0.51%  "
$ws.Range("E45").ClearFormats()

$ws.Range("B46").Value = "'PaxDollar"
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = "'1.002"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +0.16%  "
$ws.Range("E46").ClearFormats()

$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").ClearFormats()
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").ClearFormats()
$ws.Range("D47").Value = "'10.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.89%  "
$ws.Range("E47").ClearFormats()

$ws.Range("D48").Value = "'1.629"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +1.89%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").Value = "'0.06174"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -3.03%  "
$ws.Range("E49").ClearFormats()

$ws.Range("D50").Value = "'0.4446"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -2.78%  "
$ws.Range("E50").ClearFormats()

$ws.Range("D51").Value = "'1.713"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +2.26%  "
$ws.Range("E51").ClearFormats()
